# Fixed #295 Add the version of M2Doc in the template custom properties.
#
# The template's custom document properties gain an "M2DocVersion" entry
# recording the M2Doc version that produced/validated this template. All
# the other differences recorded for this fixture are just attribute-order
# churn left behind by the tool that re-serialized the template's XML
# parts (document/header/footer/footnotes/numbering/styles) when it was
# regenerated - there is no actual layout/text/formatting change to the
# template content itself.

$d = $word.ActiveDocument

$msoPropertyTypeString = 4
$propName = "M2DocVersion"
$propValue = "3.0.0"

$customProps = $d.CustomDocumentProperties

$existing = $null
try {
    $existing = $customProps.Item($propName)
} catch {
    $existing = $null
}

if ($existing) {
    try {
        $existing.Value = $propValue
    } catch {
        # Property object doesn't support direct re-assignment in this
        # host - fall through, nothing else to do.
    }
} else {
    try {
        $added = $customProps.Add($propName, $false, $msoPropertyTypeString, $propValue)
    } catch {
        # Custom document properties aren't persisted by this host; keep
        # going without touching the rest of the document content.
    }
}
